# Updates the cryptos price list (coinranking data) with the latest
# snapshot values. Column D (Price) and Column E (Volume 1h) are plain
# text cells in the source data (e.g. "29.408.07", "1.000"), so we force
# the "@" (Text) number format on column D before assigning values that
# look numeric, to keep Excel from silently re-interpreting them as
# numbers/dates and losing the original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.408.07'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.842.35'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.57'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6272'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07453'
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.83'
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.833.67'
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.973'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6772'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001027'
$ws.Range('E15').Value = '  -2.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.86'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.254'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.397.06'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '233.04'
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.335'
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.18'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.496'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1352'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.07270'
$ws.Range('E28').Value = '  +11.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.460'
$ws.Range('E29').Value = '  +3.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.481'
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.060'
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.047'
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.820'
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6981'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.573'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.993'
$ws.Range('E37').Value = '  +3.48%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.235.03'
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9460'
$ws.Range('E41').Value = '  +4.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9998'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.997.55'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.90'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.53'
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.725'
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.974'
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.921'
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1140'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('E51').Value = '  -1.36%  '
